$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename responsible person for all open non-conformities (C4:C10)
$ws.Range("C4:C10").Value = "Marisol Ornelas"

# Fill in FECHA REAL DE CIERRE (closing date) for rows 4-10
$ws.Range("E4:E10").Value = (Get-Date -Year 2016 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0).Date

# Mark rows as closed, except row 7 (id 5, "Enviar carta de aceptación") which stays "En proceso"
$ws.Range("F4:F6").Value = "Cerrada"
$ws.Range("F8:F10").Value = "Cerrada"

# Match the formatting of F4 (closed rows) across F8:F10, which previously used a distinct style
$ws.Range("F4").Copy()
$ws.Range("F8:F10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F5").Select()
